$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text with the new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.52 = 39966.67 pesos`n✅ 39966.67 pesos = 9.49 = 962.77 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the N10/O10/N12/O12 rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 105
$wsTasas.Range("O10").Value = 4196.5
$wsTasas.Range("N12").Value = 4213.5
$wsTasas.Range("O12").Value = 101.5
